$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.087.14'
$ws.Range('E2').Value = '  -6.32%  '
$ws.Range('D3').Value = '2.191.72'
$ws.Range('E3').Value = '  -7.09%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = '''240.20'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').Value = '''0.620'
$ws.Range('E6').Value = '  -7.49%  '
$ws.Range('D7').Value = '''70.25'
$ws.Range('E7').Value = '  -4.53%  '
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('D9').Value = '''0.539'
$ws.Range('E9').Value = '  -10.64%  '
$ws.Range('D10').Value = '''36.59'
$ws.Range('E10').Value = '  +7.34%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = '''0.0940'
$ws.Range('E11').Value = '  -8.07%  '
$ws.Range('B12').Value = 'OKB'
$ws.Range('C12').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D12').Value = '''57.66'
$ws.Range('E12').Value = '  -5.30%  '
$ws.Range('E13').Value = '  -4.50%  '
$ws.Range('D14').Value = '''6.55'
$ws.Range('E14').Value = '  -9.40%  '
$ws.Range('D15').Value = '2.518.96'
$ws.Range('E15').Value = '  -7.10%  '
$ws.Range('D16').Value = '''14.58'
$ws.Range('E16').Value = '  -10.04%  '
$ws.Range('D17').Value = '''0.831'
$ws.Range('E17').Value = '  -8.66%  '
$ws.Range('D18').Value = '2.185.51'
$ws.Range('E18').Value = '  -7.33%  '
$ws.Range('D19').Value = '41.002.99'
$ws.Range('E19').Value = '  -6.52%  '
$ws.Range('D20').Value = '0.0₃0939'
$ws.Range('E20').Value = '  -8.87%  '
$ws.Range('D21').Value = '''72.57'
$ws.Range('E21').Value = '  -6.56%  '
$ws.Range('D22').Value = '''6.03'
$ws.Range('E22').Value = '  -7.60%  '
$ws.Range('D23').Value = '''230.66'
$ws.Range('E23').Value = '  -8.81%  '
$ws.Range('D24').Value = '''1.99'
$ws.Range('E24').Value = '  +6.61%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('E26').Value = '  -4.96%  '
$ws.Range('D27').Value = '''2.40'
$ws.Range('E27').Value = '  -3.61%  '
$ws.Range('E28').Value = '  -5.00%  '
$ws.Range('D29').Value = '''9.69'
$ws.Range('E29').Value = '  -7.41%  '
$ws.Range('D30').Value = '''168.42'
$ws.Range('E30').Value = '  -4.40%  '
$ws.Range('D31').Value = '''20.16'
$ws.Range('E31').Value = '  -9.56%  '
$ws.Range('D32').Value = '''0.117'
$ws.Range('E32').Value = '  -8.66%  '
$ws.Range('E33').Value = '  -7.95%  '
$ws.Range('D34').Value = '''0.0699'
$ws.Range('E34').Value = '  -6.22%  '
$ws.Range('E35').Value = '  -5.10%  '
$ws.Range('E36').Value = '  -10.04%  '
$ws.Range('D37').Value = '''3.83'
$ws.Range('E37').Value = '  +1.24%  '
$ws.Range('D38').Value = '''23.04'
$ws.Range('E38').Value = '  +14.57%  '
$ws.Range('D39').Value = '''2.24'
$ws.Range('E39').Value = '  -6.81%  '
$ws.Range('E40').Value = '  -3.00%  '
$ws.Range('D41').Value = '''5.82'
$ws.Range('E41').Value = '  -11.57%  '
$ws.Range('D42').Value = '''64.19'
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('D43').Value = '''4.82'
$ws.Range('E43').Value = '  -11.28%  '
$ws.Range('D44').Value = '''8.63'
$ws.Range('E44').Value = '  -4.48%  '
$ws.Range('D45').Value = '''0.192'
$ws.Range('E45').Value = '  -4.94%  '
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').Value = '''0.0979'
$ws.Range('E47').Value = '  -7.96%  '
$ws.Range('D48').Value = '''4.47'
$ws.Range('E48').Value = '  +3.61%  '
$ws.Range('D49').Value = '''10.13'
$ws.Range('E49').Value = '  +6.20%  '
$ws.Range('E50').Value = '  -5.64%  '
$ws.Range('E51').Value = '  -6.43%  '
